# "added required experimental boolean element to valuesets"
#
# The Metadata sheet lists FHIR ValueSet properties as Property/Value rows.
# Row 7 ("Experimental") previously had no Value cell - the FHIR export now
# always emits the (boolean-as-text) "Experimental" flag, so B7 gets the
# literal text "true". The sheet's generation timestamp in row 8 ("Date")
# is also refreshed to match the regenerated export.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# B7 sits next to "Experimental" (A7) and was blank. Excel normally infers
# a literal boolean from bare "true"/"false" text, but the source file
# stores this as plain text (matching the other FHIR-boolean cell, B15's
# "BooleanType[null]") - so force text entry with a leading apostrophe,
# same as typing 'true directly into the cell.
$ws.Range("B7").Value = "'true"

# B8 sits next to "Date" (A8) - bump the generated timestamp.
$ws.Range("B8").Value = "2023-02-01T09:05:11-06:00"
